$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right total, Wrong total
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right total, Wrong total, Max display text
$ws.Range("B12").Value = 52
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "46 / 112"
